# Auto-generated Excel COM-interop script to apply the diff changes to tblStudy
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Fill in previously blank cells with "Unknown" (or "Scientists") across rows 10-48 ----
$simplePairs = @(
    @("V10", "Unknown"),
    @("AA10", "Unknown"),
    @("AC10", "Unknown"),
    @("V11", "Unknown"),
    @("AA11", "Unknown"),
    @("AC11", "Unknown"),
    @("V12", "Unknown"),
    @("AA12", "Unknown"),
    @("AC12", "Unknown"),
    @("V13", "Unknown"),
    @("AA13", "Unknown"),
    @("AC13", "Unknown"),
    @("V14", "Unknown"),
    @("AA14", "Unknown"),
    @("V15", "Unknown"),
    @("AA15", "Unknown"),
    @("AC15", "Unknown"),
    @("V16", "Unknown"),
    @("AA16", "Unknown"),
    @("AC16", "Unknown"),
    @("V17", "Unknown"),
    @("AA17", "Unknown"),
    @("AC17", "Unknown"),
    @("V18", "Unknown"),
    @("AA18", "Unknown"),
    @("AC18", "Unknown"),
    @("V19", "Unknown"),
    @("AA19", "Unknown"),
    @("AC19", "Unknown"),
    @("V20", "Unknown"),
    @("AA20", "Unknown"),
    @("AC20", "Unknown"),
    @("V22", "Unknown"),
    @("AA22", "Unknown"),
    @("AC22", "Unknown"),
    @("V23", "Unknown"),
    @("AA23", "Unknown"),
    @("AC23", "Unknown"),
    @("V25", "Unknown"),
    @("AA25", "Unknown"),
    @("AC25", "Unknown"),
    @("V26", "Unknown"),
    @("AC26", "Unknown"),
    @("V27", "Unknown"),
    @("AA27", "Unknown"),
    @("W28", "Unknown"),
    @("V29", "Unknown"),
    @("AA29", "Unknown"),
    @("AC29", "Unknown"),
    @("AA30", "Unknown"),
    @("AC30", "Unknown"),
    @("V31", "Unknown"),
    @("AA31", "Unknown"),
    @("AC31", "Unknown"),
    @("V32", "Unknown"),
    @("AA32", "Unknown"),
    @("AC32", "Unknown"),
    @("V33", "Unknown"),
    @("AA33", "Unknown"),
    @("AC33", "Unknown"),
    @("V34", "Unknown"),
    @("AA34", "Unknown"),
    @("AC34", "Unknown"),
    @("V35", "Unknown"),
    @("AA35", "Unknown"),
    @("AC35", "Unknown"),
    @("V36", "Unknown"),
    @("Y36", "Unknown"),
    @("Z36", "Scientists"),
    @("AA36", "Unknown"),
    @("AC36", "Unknown"),
    @("AD36", "Unknown"),
    @("V37", "Unknown"),
    @("AA37", "Unknown"),
    @("AC37", "Unknown"),
    @("V38", "Unknown"),
    @("AA38", "Unknown"),
    @("AC38", "Unknown"),
    @("V39", "Unknown"),
    @("AA39", "Unknown"),
    @("V40", "Unknown"),
    @("AA40", "Unknown"),
    @("AC40", "Unknown"),
    @("V41", "Unknown"),
    @("AA41", "Unknown"),
    @("AC41", "Unknown"),
    @("V42", "Unknown"),
    @("AA42", "Unknown"),
    @("AC42", "Unknown"),
    @("V43", "Unknown"),
    @("AA43", "Unknown"),
    @("AC43", "Unknown"),
    @("V45", "Unknown"),
    @("AA45", "Unknown"),
    @("AC45", "Unknown"),
    @("V46", "Unknown"),
    @("AA46", "Unknown"),
    @("AB46", "Unknown"),
    @("AC46", "Unknown"),
    @("AD46", "Unknown"),
    @("V47", "Unknown"),
    @("AA47", "Unknown"),
    @("AC47", "Unknown"),
    @("AD47", "Unknown"),
    @("AA48", "Unknown"),
    @("AC48", "Unknown")
)
foreach ($pair in $simplePairs) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- 2. Species (column D) text edits - punctuation / wording clean-up ----
$dChanges = @(
    @("D10", "red king crab (Paralithodes camtschaticus)"),
    @("D11", "common coral trout (Plectropomus leopardus), red throat emperor (Lethrinus miniatus)"),
    @("D12", "tiger flathead (Neoplatycephalus richardsoni), jackass morwong (Nemadactylus macropterus), school whiting (Sillago flindersi)"),
    @("D16", "Tiger prawn (Penaeus esculentus), Tiger prawn (P. semisulcatus), Endeavour prawn (Metapenaeus endeavouri), Endeavour prawn (Metapenaeus ensis), Northern Australia prawn ecosystem"),
    @("D17", "Megrim (Lepidorhombus whiffiagonis), Hake (Merluccius merluccius), Black anglerfish (Lophius budegassa), White anglerfish (Lophius piscatorius), Western Horse mackeral (Trachurus trachurus), Mackeral (Scomber scombrus), Blue whiting (Micromesistius poutassou), Rays (Leucoraja naevus), Inshore squids (Loliginidae, Seabass (Dicentrarchus labrax), Cuttlefish (Sepiidae), bobtail squids (Sepiolidae), Red mullet (Mullus surmuletus)"),
    @("D18", "Tiger prawn (Penaeus esculentus), Tiger prawn (Penaeus semisulcatus)"),
    @("D31", "abalone (Haliotis rubra), abalone (Haliotis laevigata)")
)
foreach ($pair in $dChanges) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# ---- 3. IncludeInPublication flag fix for row 45 (now included) ----
$ws.Range("AI45").Value = $true

# ---- 4. Append five new study rows (49-53) ----
$newRows = @(
    @{ "row" = 49; "cells" = @(
        @("A", 51),
        @("B", "A'mar et al"),
        @("C", "2009"),
        @("D", "walleye pollock (Theragra chalcogramma)"),
        @("E", "Gulf of Alaska"),
        @("F", 57),
        @("G", -144),
        @("H", "Gulf of Alaska walleye pollock fishery"),
        @("I", "Climate Change"),
        @("J", $false),
        @("K", $false),
        @("L", $false),
        @("M", $false),
        @("N", $false),
        @("O", $false),
        @("P", $false),
        @("Q", $false),
        @("R", $false),
        @("S", $false),
        @("T", "My definition: `"How well do the current management strategy, and some candidate management strategies, perform given shifts in climate regime, as realized in fisheries productivity in the Gulf of Alaska walleye pollock fishery?`""),
        @("U", "Not discussed, but apparently they were supplied by the scientists."),
        @("V", "Unknown"),
        @("W", "Unknown"),
        @("X", ""),
        @("Y", "Scientists"),
        @("Z", "Scientists"),
        @("AA", "Unknown"),
        @("AB", "Scientists"),
        @("AC", "Unknown"),
        @("AD", "Scientists"),
        @("AE", "Simulation modeling"),
        @("AF", "A’mar, Z. Teresa, André E. Punt, and Martin W. Dorn. “The Impact of Regime Shifts on the Performance of Management Strategies for the Gulf of Alaska Walleye Pollock (Theragra Chalcogramma) Fishery.” Canadian Journal of Fisheries and Aquatic Sciences 66, no. 12 (December 2009): 2222–42."),
        @("AG", "doi.org/10.1139/F09-142"),
        @("AH", "This was a simulation study primarily."),
        @("AI", $true)
    ) },
    @{ "row" = 50; "cells" = @(
        @("A", 52),
        @("B", "Dorner et al"),
        @("C", "2009"),
        @("D", "Pacific salmon (Onchorhyncus spp.)"),
        @("E", "North Pacific"),
        @("F", 51),
        @("G", -129),
        @("H", "Pacific Salmon fsihery"),
        @("I", "Climate Change"),
        @("J", $false),
        @("K", $false),
        @("L", $false),
        @("M", $false),
        @("N", $false),
        @("O", $false),
        @("P", $false),
        @("Q", $false),
        @("R", $false),
        @("S", $false),
        @("T", "The authors state: `"the purpose of our research was to evaluate the relative performance of several combinations of harvest policies and stock assessment/forecasting models, including hierarchical models and models that make use of environmental covariates, in the presence of uncertainties about future climatic conditions and outcome uncertainty`""),
        @("U", "Seemingly supplied by the scientist authoring the MSE"),
        @("V", "Unknown"),
        @("W", "Unknown"),
        @("X", ""),
        @("Y", "Unknown"),
        @("Z", "Scientists"),
        @("AA", "Unknown"),
        @("AB", "Scientists"),
        @("AC", "Unknown"),
        @("AD", "Scientists"),
        @("AE", "Simulation modeling"),
        @("AF", "Dorner, Brigitte, Randall M. Peterman, and Zhenming Su. “Evaluation of Performance of Alternative Management Models of Pacific Salmon (Oncorhynchus Spp.) in the Presence of Climatic Change and Outcome Uncertainty Using Monte Carlo Simulations.” Canadian Journal of Fisheries and Aquatic Sciences 66, no. 12 (December 2009): 2199–2221."),
        @("AG", "10.1139/F09-144"),
        @("AH", "This is a simulation study priamrily focused on alternative stock assessment methods"),
        @("AI", $true)
    ) },
    @{ "row" = 51; "cells" = @(
        @("A", 53),
        @("B", "Haltuch et al"),
        @("C", "2019"),
        @("D", "Sablefish (Anoplopoma fimbria)"),
        @("E", "North Pacific Ocean off US West Coast"),
        @("F", 41),
        @("G", -126),
        @("H", "West Coast sablefish fishery"),
        @("I", "Climate Change"),
        @("J", $false),
        @("K", $false),
        @("L", $false),
        @("M", $false),
        @("N", $false),
        @("O", $false),
        @("P", $false),
        @("Q", $false),
        @("R", $false),
        @("S", $false),
        @("T", "`"how resilient is the sablefish stock is to current fishery harvest control rules (HCRs) given climate change and variability and considering alternatives that might be more responsive to long-term directional changes in the productivity of fish stocks?`""),
        @("U", "Seemingly supplied by scientists"),
        @("V", "Unknown"),
        @("W", "Unknown"),
        @("X", ""),
        @("Y", "Unknown"),
        @("Z", "Scientists"),
        @("AA", "Unknown"),
        @("AB", "Scientists"),
        @("AC", "Unknown"),
        @("AD", "Scientists"),
        @("AE", "Simulation modeling"),
        @("AF", "Haltuch, Melissa A, Z Teresa A’mar, Nicholas A Bond, and Juan L Valero. “Assessing the Effects of Climate Change on US West Coast Sablefish Productivity and on the Performance of Alternative Management Strategies.” Edited by Jörn Schmidt. ICES Journal of Marine Science 76, no. 6 (December 1, 2019): 1524–42."),
        @("AG", "10.1093/icesjms/fsz029"),
        @("AH", "Seemingly a simulation study without any connection to the management process."),
        @("AI", $true)
    ) },
    @{ "row" = 52; "cells" = @(
        @("A", 54),
        @("B", "Merino et al"),
        @("C", "2019"),
        @("D", "North Atlantic Albacore (Thunnus alalunga)"),
        @("E", "North Atlantic Ocean"),
        @("F", 38),
        @("G", -39),
        @("H", "North Atlantic Albacore Fishery"),
        @("I", "Climate Change"),
        @("J", $false),
        @("K", $false),
        @("L", $false),
        @("M", $false),
        @("N", $false),
        @("O", $false),
        @("P", $false),
        @("Q", $false),
        @("R", $false),
        @("S", $false),
        @("T", "My problem statement based on the available documentation, `"Is the HCR adopted for North Atlantic albacore robust to a range of climate change impacts`""),
        @("U", "Seemingly they were provided by the scientists conducting the evaluation based on established policy from the management body, but the methodology was not explictly documented."),
        @("V", "Mental Analysis"),
        @("W", "Mental Analysis"),
        @("X", ""),
        @("Y", "Unknown"),
        @("Z", "Scientists"),
        @("AA", "Unknown"),
        @("AB", "Management, Scientists"),
        @("AC", "Unknown"),
        @("AD", "Scientists"),
        @("AE", "Simulation modeling"),
        @("AF", "Merino, Gorka, Haritz Arrizabalaga, Igor Arregui, Josu Santiago, Hilario Murua, Agurtzane Urtizberea, Eider Andonegi, Paul De Bruyn, and Laurence T. Kell. “Adaptation of North Atlantic Albacore Fishery to Climate Change: Yet Another Potential Benefit of Harvest Control Rules.” Frontiers in Marine Science 6 (October 10, 2019)"),
        @("AG", "10.3389/fmars.2019.00620"),
        @("AH", "Seemingly a simulation study without a decision making process. No trade-off as evaluated as there are no alternative management proceedures to compare."),
        @("AI", $true)
    ) },
    @{ "row" = 53; "cells" = @(
        @("A", 55),
        @("B", "Castillo-Jordán et al"),
        @("C", "2019"),
        @("D", "Patagonian Grenadier (Macruronus Magellanicus)"),
        @("E", "Chile"),
        @("F", -44),
        @("G", -76),
        @("H", "Chile Patagonian Grenadier fishery"),
        @("I", "Climate Change"),
        @("J", $false),
        @("K", $false),
        @("L", $false),
        @("M", $false),
        @("N", $false),
        @("O", $false),
        @("P", $false),
        @("Q", $false),
        @("R", $false),
        @("S", $false),
        @("T", "My slight adaptation of the study objective is: `"how would a regime shift in recruitment for Patagonian grenadier off Chile impact the success of management decisions and the sustainability of the fishery?`""),
        @("U", "Not addressed"),
        @("V", "Unknown"),
        @("W", "Mental Analysis"),
        @("X", ""),
        @("Y", "Unknown"),
        @("Z", "Scientists"),
        @("AA", "Unknown"),
        @("AB", "Scientists"),
        @("AC", "Unknown"),
        @("AD", "Scientists"),
        @("AE", "Simulation modeling"),
        @("AF", "Castillo-Jordán, Claudio, Sally E. Wayte, Geoffrey N. Tuck, Sean Tracey, Stewart D. Frusher, and André E. Punt. “Implications of a Climate-Induced Recruitment Shift in the Stock Assessment of Patagonian Grenadier (Macruronus Magellanicus) in Chile.” Fisheries Research 212 (April 2019): 114–22."),
        @("AG", "10.1016/j.fishres.2018.12.019"),
        @("AH", "While this is referred to as a climate change study, there isn't a projection of climate change explicitly. What is projected is a possible regime shift in the environment. In terms of the MSE this is a simulation study with an assoiciated decision process seemingly."),
        @("AI", $true)
    ) }
)

# Columns that must always be written as text, even when the content looks numeric
$textForceCols = @("C")

# Columns that are right-aligned in the existing table (matches styles used for ID/lat/long/boolean columns)
$rightAlignCols = @("A", "AI", "F", "G", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S")

foreach ($rowDef in $newRows) {
    $r = $rowDef["row"]
    foreach ($cellPair in $rowDef["cells"]) {
        $colLetter = $cellPair[0]
        $cellValue = $cellPair[1]
        $target = $ws.Range("$colLetter$r")
        if ($textForceCols -contains $colLetter) {
            $target.NumberFormat = "@"
        }
        $target.Value = $cellValue
        $target.WrapText = $true
        $target.VerticalAlignment = -4108
        if ($rightAlignCols -contains $colLetter) {
            $target.HorizontalAlignment = -4152
        }
    }
}

